# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q1" sheet, placed right before "总计" ---
$totalSheetRef = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetRef)
$newSheet.Name = "2022-Q1"

# NOTE: the "Before" reference used above becomes stale once the sheet
# collection shifts around - always re-look-up "总计" by name afterwards.
$totalSheet = $wb.Worksheets.Item("总计")

# Reuse the header/column styling (bold header row + bold index column)
# from an existing per-quarter sheet so the new sheet matches the
# established look (style index 2 = bold, centered, bordered).
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$srcSheet.Range("A1:H8").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Fill in the header row ---
$newSheet.Cells.Item(1,2).Value2 = "基金代码"
$newSheet.Cells.Item(1,3).Value2 = "基金名称"
$newSheet.Cells.Item(1,4).Value2 = "基金规模"
$newSheet.Cells.Item(1,5).Value2 = "股票总仓位"
$newSheet.Cells.Item(1,6).Value2 = "仓位占比"
$newSheet.Cells.Item(1,7).Value2 = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value2 = "仓位排名"

# --- 3. Fill in the fund holdings data rows ---
# Columns D, E, F, G are text-formatted numeric-looking figures (kept as
# text so values like "34.50" don't lose their trailing zero), B is a
# fund code that must keep leading zeros, H is a genuine integer rank.
function Set-FundRow($row, $idx, $code, $name, $scale, $pos, $ratio, $value, $rank) {
    $newSheet.Cells.Item($row,1).Value2 = $idx
    $newSheet.Cells.Item($row,2).Value2 = "'" + $code
    $newSheet.Cells.Item($row,3).Value2 = $name
    $newSheet.Cells.Item($row,4).Value2 = "'" + $scale
    $newSheet.Cells.Item($row,5).Value2 = "'" + $pos
    $newSheet.Cells.Item($row,6).Value2 = "'" + $ratio
    $newSheet.Cells.Item($row,7).Value2 = "'" + $value
    $newSheet.Cells.Item($row,8).Value2 = $rank
}

Set-FundRow 2 0 "009010" "华夏兴阳一年持有期混合" "34.50" "90.35" "2.34" "0.8073" 10
Set-FundRow 3 1 "008850" "景顺长城价值稳进三年定期开放灵活配置混合" "17.06" "69.71" "2.21" "0.3770" 7
Set-FundRow 4 2 "008715" "景顺长城价值驱动一年持有期灵活配置混合型证券投资基金" "16.83" "62.03" "2.20" "0.3703" 7
Set-FundRow 5 3 "000979" "景顺长城沪港深精选股票" "16.46" "82.61" "2.12" "0.3490" 10
Set-FundRow 6 4 "009098" "景顺长城价值领航两年持有期混合" "11.67" "75.58" "2.33" "0.2719" 6
Set-FundRow 7 5 "008060" "景顺长城价值边际灵活配置混合" "4.93" "80.78" "2.23" "0.1099" 9
Set-FundRow 8 6 "006205" "汇添富沪港深优势精选定期开放混合" "0.40" "93.67" "4.21" "0.0168" 9

# --- 4. Insert a new top data row in "总计" for the 2022-Q1 summary,
#        shifting the existing quarters down by one row ---
$totalSheet.Rows.Item(2).Insert()

# The row below (now row 3, the old 2021-Q4 row) still carries the
# correctly-established data-row styling (bold/bordered index cell,
# plain value cells) - copy just that 4-cell format onto the new blank
# row instead of Font-toggling it cell by cell.
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Cells.Item(2,1).Value2 = 0
$totalSheet.Cells.Item(2,2).Value2 = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value2 = 7
$totalSheet.Cells.Item(2,4).Value2 = 2.3

# Renumber the shifted index column (A) to stay 0,1,2,3,4
$totalSheet.Cells.Item(3,1).Value2 = 1
$totalSheet.Cells.Item(4,1).Value2 = 2
$totalSheet.Cells.Item(5,1).Value2 = 3
$totalSheet.Cells.Item(6,1).Value2 = 4

# --- 5. Restore the originally-active sheet/selection (index 1, "2021-Q1")
#        so this edit doesn't also change the saved view state. ---
$wb.Worksheets.Item(1).Activate()
$null = $wb.Worksheets.Item(1).Range("A1").Select()
